$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe everything currently used on the sheet (contents + formatting)
$ws.Cells.Clear()

# Write the single remaining row back at row 3 (A3 holds the text, B3:D3 share its style but stay empty)
$ws.Range("A3").Value = "tes ke dua x"
$ws.Range("A3:D3").Interior.Color = 65535   # 0x00FFFF00 -> RGB(255,255,0) yellow

$ws.Range("F2").Select()
